# Update cryptos list with refreshed prices / volume percentages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.758.24'
$ws.Range('E2').Value = '  +1.41%  '
$ws.Range('D3').Value = '1.877.34'
$ws.Range('E3').Value = '  +1.08%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '331.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.54%  '
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4724'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3952'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.01'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08072'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.031'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.17'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.73%  '
$ws.Range('D13').Value = '1.877.57'
$ws.Range('E13').Value = '  +1.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.971'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.146'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.006'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('E17').Value = '  +1.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '87.11'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06658'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').Value = '27.771.16'
$ws.Range('E22').Value = '  +1.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.529'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.308'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.27%  '
$ws.Range('D26').Value = '2.106.51'
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '158.89'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.110'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.588'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '122.41'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9860'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09553'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.454'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.589'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.350'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.61%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06118'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02259'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.228'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.174'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6038'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1906'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.68%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5729'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.00%  '
$ws.Range('B45').Value = 'WEMIXTOKEN'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.252'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.25'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.950'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.382'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06899'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '115.04'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.073'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.90%  '
